$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the updated price/volume cells stay plain text (matching the
# workbook author's original inline-string storage) instead of being
# auto-converted to numbers/percentages by Excel.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '300.81'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '-4.36%'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '35.17'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '-0.36%'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.047'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '-1.70%'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.07978'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '-1.84%'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.901'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '-10.05%'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '7.806'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '-1.90%'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '4.050'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '-2.54%'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.931'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '0.79%'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9215'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '-0.77%'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.1278'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '26.86%'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.1849'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '-1.32%'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.09740'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '6.73%'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.03568'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '-1.15%'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.09852'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.001392'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '-2.85%'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.005765'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '1.73%'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.3398'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '-0.43%'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.1292'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '-3.69%'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.059'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '-0.73%'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '8.43%'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.04506'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '-1.52%'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.001213'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '-2.76%'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.004782'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '1.78%'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0003002'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '-33.41%'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01874'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '-4.48%'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.04686'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '-3.28%'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.007486'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '-3.03%'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.01021'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '30.03%'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1324'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '-4.80%'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '0.00%'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.01067'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '-9.57%'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00006249'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '-6.18%'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.00000000750'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '-0.13%'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '63.17%'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '-12.55%'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.00002101'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '-0.13%'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0002001'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '-0.13%'
